$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFCB1")

# Insert a new column at E ("chl hv"), shifting the existing E:P data right to F:Q
$ws.Columns("E:E").Insert()
# The insert stamps inherited formatting into every existing row of the new
# column; clear those placeholder cells since column E has no data in rows 2:13
$ws.Range("E2:E13").Clear()

# --- New row 14: partial entry noting the new laser after a lightning strike ---
$ws.Range("B14").Value = 42389
$ws.Range("B14").NumberFormat = "m/d/yy"
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("D14").Value = "4.5-5V"
$ws.Range("E1").Value = "chl hv"
$ws.Range("E14").Value = 0.7
$ws.Range("E14").HorizontalAlignment = -4108
$ws.Range("E14").VerticalAlignment = -4108
$ws.Range("N14").Value = "new laser after lightening strike"

# --- New row 15: IFCB1_036_165039 ---
$ws.Range("A15").Value = "IFCB1_036_165039"
$ws.Range("B15").Value = 42405
$ws.Range("B15").NumberFormat = "m/d/yy"
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("B15").VerticalAlignment = -4108
$ws.Range("C15").Value = 165039
$ws.Range("D15").Value = "4V"
$ws.Range("E15").Value = 0.7
$ws.Range("E15").HorizontalAlignment = -4108
$ws.Range("E15").VerticalAlignment = -4108
$ws.Range("F15").Value = 2.8
$ws.Range("G15").Value = 5.8
$ws.Range("H15").Value = 29.8
$ws.Range("M15").Value = "0.3-0.6"
$ws.Range("I15").Value = "2.5-3.2"
$ws.Range("J15").Value = 0.42
$ws.Range("K15").Value = 9.8
$ws.Range("L15").Value = 35.7

# --- New row 16: IFCB1_036_172753 ---
$ws.Range("A16").Value = "IFCB1_036_172753"
$ws.Range("B16").Value = 42405
$ws.Range("B16").NumberFormat = "m/d/yy"
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("C16").Value = 172753
$ws.Range("E16").Value = 0.7
$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("E16").VerticalAlignment = -4108
$ws.Range("F16").Value = 2.8
$ws.Range("G16").Value = 5.8
$ws.Range("H16").Value = 26.9
$ws.Range("I16").Value = "2.5-3.1"
$ws.Range("J16").Value = 0.41
$ws.Range("K16").Value = 12.2
$ws.Range("L16").Value = 38.7
$ws.Range("M16").Value = "0.3-0.6"

# --- comments2 column for rows 15/16, written last ---
$ws.Range("N15").Value = "use all signals, realign who camera/PMT stack,pump1"
$ws.Range("N16").Value = "use all signals, realign who camera/PMT stack,pump2"

# Make IFCB1 the active sheet/tab with N17 selected (was IFCB5 before the edit)
$ws.Activate()
$ws.Range("N17").Select()
